$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and two pairs of
# re-ordered rows) as captured by the GitHub Actions scraper run.
# Each cell is forced to Text format before the write (and the format
# is cleared again afterwards) so Excel stores the numeric-looking
# strings (e.g. "1.000", "0.6834") verbatim instead of silently
# coercing them into numbers.

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = '@'
$cell.Value = '29.151.90'
$cell.ClearFormats()
$cell = $ws.Range("E2")
$cell.NumberFormat = '@'
$cell.Value = '  +0.72%  '
$cell.ClearFormats()
# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = '@'
$cell.Value = '1.833.77'
$cell.ClearFormats()
$cell = $ws.Range("E3")
$cell.NumberFormat = '@'
$cell.Value = '  -0.07%  '
$cell.ClearFormats()
# Row 4
$cell = $ws.Range("D4")
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.ClearFormats()
$cell = $ws.Range("E4")
$cell.NumberFormat = '@'
$cell.Value = '  +0.14%  '
$cell.ClearFormats()
# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = '@'
$cell.Value = '240.09'
$cell.ClearFormats()
$cell = $ws.Range("E5")
$cell.NumberFormat = '@'
$cell.Value = '  -2.13%  '
$cell.ClearFormats()
# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = '@'
$cell.Value = '0.6834'
$cell.ClearFormats()
$cell = $ws.Range("E6")
$cell.NumberFormat = '@'
$cell.Value = '  -1.37%  '
$cell.ClearFormats()
# Row 7
$cell = $ws.Range("E7")
$cell.NumberFormat = '@'
$cell.Value = '  +0.11%  '
$cell.ClearFormats()
# Row 8
$cell = $ws.Range("E8")
$cell.NumberFormat = '@'
$cell.Value = '  -1.10%  '
$cell.ClearFormats()
# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = '@'
$cell.Value = '0.07457'
$cell.ClearFormats()
$cell = $ws.Range("E9")
$cell.NumberFormat = '@'
$cell.Value = '  -3.06%  '
$cell.ClearFormats()
# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = '@'
$cell.Value = '23.07'
$cell.ClearFormats()
$cell = $ws.Range("E10")
$cell.NumberFormat = '@'
$cell.Value = '  -1.33%  '
$cell.ClearFormats()
# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = '@'
$cell.Value = '0.07651'
$cell.ClearFormats()
$cell = $ws.Range("E11")
$cell.NumberFormat = '@'
$cell.Value = '  -2.04%  '
$cell.ClearFormats()
# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = '@'
$cell.Value = '1.842.34'
$cell.ClearFormats()
$cell = $ws.Range("E12")
$cell.NumberFormat = '@'
$cell.Value = '  +0.50%  '
$cell.ClearFormats()
# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = '@'
$cell.Value = '5.055'
$cell.ClearFormats()
$cell = $ws.Range("E13")
$cell.NumberFormat = '@'
$cell.Value = '  -1.01%  '
$cell.ClearFormats()
# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = '@'
$cell.Value = '0.6824'
$cell.ClearFormats()
$cell = $ws.Range("E14")
$cell.NumberFormat = '@'
$cell.Value = '  -0.06%  '
$cell.ClearFormats()
# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = '@'
$cell.Value = '87.06'
$cell.ClearFormats()
$cell = $ws.Range("E15")
$cell.NumberFormat = '@'
$cell.Value = '  -6.78%  '
$cell.ClearFormats()
# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = '@'
$cell.Value = '6.218'
$cell.ClearFormats()
$cell = $ws.Range("E16")
$cell.NumberFormat = '@'
$cell.Value = '  -5.62%  '
$cell.ClearFormats()
# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = '@'
$cell.Value = '29.163.68'
$cell.ClearFormats()
$cell = $ws.Range("E17")
$cell.NumberFormat = '@'
$cell.Value = '  +0.80%  '
$cell.ClearFormats()
# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = '@'
$cell.Value = '0.000008175'
$cell.ClearFormats()
$cell = $ws.Range("E18")
$cell.NumberFormat = '@'
$cell.Value = '  -1.13%  '
$cell.ClearFormats()
# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = '@'
$cell.Value = '2.080.51'
$cell.ClearFormats()
$cell = $ws.Range("E19")
$cell.NumberFormat = '@'
$cell.Value = '  +0.37%  '
$cell.ClearFormats()
# Row 20
$cell = $ws.Range("B20")
$cell.NumberFormat = '@'
$cell.Value = 'Avalanche'
$cell.ClearFormats()
$cell = $ws.Range("C20")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell.ClearFormats()
$cell = $ws.Range("D20")
$cell.NumberFormat = '@'
$cell.Value = '12.53'
$cell.ClearFormats()
$cell = $ws.Range("E20")
$cell.NumberFormat = '@'
$cell.Value = '  -1.36%  '
$cell.ClearFormats()
# Row 21
$cell = $ws.Range("B21")
$cell.NumberFormat = '@'
$cell.Value = 'BitcoinCash'
$cell.ClearFormats()
$cell = $ws.Range("C21")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell.ClearFormats()
$cell = $ws.Range("D21")
$cell.NumberFormat = '@'
$cell.Value = '226.38'
$cell.ClearFormats()
$cell = $ws.Range("E21")
$cell.NumberFormat = '@'
$cell.Value = '  -5.98%  '
$cell.ClearFormats()
# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.ClearFormats()
$cell = $ws.Range("E22")
$cell.NumberFormat = '@'
$cell.Value = '  +0.05%  '
$cell.ClearFormats()
# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = '@'
$cell.Value = '7.400'
$cell.ClearFormats()
$cell = $ws.Range("E23")
$cell.NumberFormat = '@'
$cell.Value = '  -0.85%  '
$cell.ClearFormats()
# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell.ClearFormats()
$cell = $ws.Range("E24")
$cell.NumberFormat = '@'
$cell.Value = '  +0.11%  '
$cell.ClearFormats()
# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = '@'
$cell.Value = '0.1456'
$cell.ClearFormats()
$cell = $ws.Range("E25")
$cell.NumberFormat = '@'
$cell.Value = '  -3.15%  '
$cell.ClearFormats()
# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = '@'
$cell.Value = '159.65'
$cell.ClearFormats()
$cell = $ws.Range("E26")
$cell.NumberFormat = '@'
$cell.Value = '  +0.27%  '
$cell.ClearFormats()
# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = '@'
$cell.Value = '8.768'
$cell.ClearFormats()
$cell = $ws.Range("E27")
$cell.NumberFormat = '@'
$cell.Value = '  +0.09%  '
$cell.ClearFormats()
# Row 28
$cell = $ws.Range("E28")
$cell.NumberFormat = '@'
$cell.Value = '  -0.76%  '
$cell.ClearFormats()
# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = '@'
$cell.Value = '1.504'
$cell.ClearFormats()
$cell = $ws.Range("E29")
$cell.NumberFormat = '@'
$cell.Value = '  -2.51%  '
$cell.ClearFormats()
# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = '@'
$cell.Value = '4.259'
$cell.ClearFormats()
$cell = $ws.Range("E30")
$cell.NumberFormat = '@'
$cell.Value = '  +0.97%  '
$cell.ClearFormats()
# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = '@'
$cell.Value = '4.143'
$cell.ClearFormats()
$cell = $ws.Range("E31")
$cell.NumberFormat = '@'
$cell.Value = '  -0.15%  '
$cell.ClearFormats()
# Row 32
$cell = $ws.Range("E32")
$cell.NumberFormat = '@'
$cell.Value = '  +1.12%  '
$cell.ClearFormats()
# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = '@'
$cell.Value = '0.05140'
$cell.ClearFormats()
$cell = $ws.Range("E33")
$cell.NumberFormat = '@'
$cell.Value = '  +0.62%  '
$cell.ClearFormats()
# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = '@'
$cell.Value = '0.7678'
$cell.ClearFormats()
$cell = $ws.Range("E34")
$cell.NumberFormat = '@'
$cell.Value = '  -1.19%  '
$cell.ClearFormats()
# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = '@'
$cell.Value = '1.844'
$cell.ClearFormats()
$cell = $ws.Range("E35")
$cell.NumberFormat = '@'
$cell.Value = '  -0.77%  '
$cell.ClearFormats()
# Row 36
$cell = $ws.Range("E36")
$cell.NumberFormat = '@'
$cell.Value = '  -1.37%  '
$cell.ClearFormats()
# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = '@'
$cell.Value = '2.673'
$cell.ClearFormats()
$cell = $ws.Range("E37")
$cell.NumberFormat = '@'
$cell.Value = '  -0.78%  '
$cell.ClearFormats()
# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = '@'
$cell.Value = '1.305.14'
$cell.ClearFormats()
$cell = $ws.Range("E38")
$cell.NumberFormat = '@'
$cell.Value = '  +0.91%  '
$cell.ClearFormats()
# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = '@'
$cell.Value = '0.01835'
$cell.ClearFormats()
$cell = $ws.Range("E39")
$cell.NumberFormat = '@'
$cell.Value = '  -1.29%  '
$cell.ClearFormats()
# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = '@'
$cell.Value = '2.701'
$cell.ClearFormats()
$cell = $ws.Range("E40")
$cell.NumberFormat = '@'
$cell.Value = '  +0.08%  '
$cell.ClearFormats()
# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = '@'
$cell.Value = '0.9353'
$cell.ClearFormats()
$cell = $ws.Range("E41")
$cell.NumberFormat = '@'
$cell.Value = '  -1.93%  '
$cell.ClearFormats()
# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = '@'
$cell.Value = '5.829'
$cell.ClearFormats()
$cell = $ws.Range("E42")
$cell.NumberFormat = '@'
$cell.Value = '  -5.89%  '
$cell.ClearFormats()
# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = '@'
$cell.Value = '104.25'
$cell.ClearFormats()
$cell = $ws.Range("E43")
$cell.NumberFormat = '@'
$cell.Value = '  -2.35%  '
$cell.ClearFormats()
# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = '@'
$cell.Value = '0.9997'
$cell.ClearFormats()
$cell = $ws.Range("E44")
$cell.NumberFormat = '@'
$cell.Value = '  +0.05%  '
$cell.ClearFormats()
# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = '@'
$cell.Value = '65.61'
$cell.ClearFormats()
$cell = $ws.Range("E45")
$cell.NumberFormat = '@'
$cell.Value = '  +2.65%  '
$cell.ClearFormats()
# Row 46
$cell = $ws.Range("B46")
$cell.NumberFormat = '@'
$cell.Value = 'RocketPoolETH'
$cell.ClearFormats()
$cell = $ws.Range("C46")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$cell.ClearFormats()
$cell = $ws.Range("D46")
$cell.NumberFormat = '@'
$cell.Value = '1.980.62'
$cell.ClearFormats()
$cell = $ws.Range("E46")
$cell.NumberFormat = '@'
$cell.Value = '  +0.32%  '
$cell.ClearFormats()
# Row 47
$cell = $ws.Range("B47")
$cell.NumberFormat = '@'
$cell.Value = 'Mantle'
$cell.ClearFormats()
$cell = $ws.Range("C47")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$cell.ClearFormats()
$cell = $ws.Range("D47")
$cell.NumberFormat = '@'
$cell.Value = '0.5196'
$cell.ClearFormats()
$cell = $ws.Range("E47")
$cell.NumberFormat = '@'
$cell.Value = '  +0.60%  '
$cell.ClearFormats()
# Row 48
$cell = $ws.Range("B48")
$cell.NumberFormat = '@'
$cell.Value = 'EnergySwap'
$cell.ClearFormats()
$cell = $ws.Range("C48")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell.ClearFormats()
$cell = $ws.Range("D48")
$cell.NumberFormat = '@'
$cell.Value = '9.565'
$cell.ClearFormats()
$cell = $ws.Range("E48")
$cell.NumberFormat = '@'
$cell.Value = '  -1.39%  '
$cell.ClearFormats()
# Row 49
$cell = $ws.Range("B49")
$cell.NumberFormat = '@'
$cell.Value = 'RenderToken'
$cell.ClearFormats()
$cell = $ws.Range("C49")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell.ClearFormats()
$cell = $ws.Range("D49")
$cell.NumberFormat = '@'
$cell.Value = '1.769'
$cell.ClearFormats()
$cell = $ws.Range("E49")
$cell.NumberFormat = '@'
$cell.Value = '  +0.66%  '
$cell.ClearFormats()
# Row 50
$cell = $ws.Range("B50")
$cell.NumberFormat = '@'
$cell.Value = 'BabyDogeCoin'
$cell.ClearFormats()
$cell = $ws.Range("C50")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$cell.ClearFormats()
$cell = $ws.Range("D50")
$cell.NumberFormat = '@'
$cell.Value = '0.00000000120'
$cell.ClearFormats()
$cell = $ws.Range("E50")
$cell.NumberFormat = '@'
$cell.Value = '  -2.57%  '
$cell.ClearFormats()
# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = '@'
$cell.Value = '0.07391'
$cell.ClearFormats()
$cell = $ws.Range("E51")
$cell.NumberFormat = '@'
$cell.Value = '  +21.02%  '
$cell.ClearFormats()
